# Correct typos & update offloading fig
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in the "quantize_percentile" description (row 6, column C):
# add missing periods at the end of the 2nd and 3rd lines.
$ws.Range("C6").Value = "Percentile used for the quantization method ""Percentile"" and ""MaxPercentile"".`nThis should be between 0 and 1. (Ex. 0.999, 0.9999).`nDefaults to 0.9999."

# Update the active selection to C10 (reflects the updated offloading figure view)
$ws.Range("C10").Select()
